$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the old row 39 ("MailBoxes" section), which
# pushes the MailBoxes / Regex sections down by two rows (old 39-53 -> 41-55).
$ws.Rows.Item(39).Insert() | Out-Null
$ws.Rows.Item(39).Insert() | Out-Null

# Row heights for the two new rows match the other wrapped-text rows (30pt).
$ws.Rows.Item(39).RowHeight = 30
$ws.Rows.Item(40).RowHeight = 30

# Row 39: InProgressWindowFolder
$ws.Range("A39").Value = "InProgressWindowFolder"
$ws.Range("B39").Value = "\\cavmfil001\Common\SinglePaymentScheme\Cross Compliance {0}\9)Processing Folders\Robot\{1}\In Progress\"
$ws.Range("C39").Value = "in progress path for all crf types"

# Row 40: ReadyToBeReceiptedFolder
$ws.Range("A40").Value = "ReadyToBeReceiptedFolder"
$ws.Range("B40").Value = "\\cavmfil001\Common\SinglePaymentScheme\Cross Compliance {0}\9)Processing Folders\Robot\{1}\Ready To Be Receipted\"
$ws.Range("C40").Value = "ready to be receipted for all crf types - use for deployment 1 only"

# Grow the Table1 ListObject (and its AutoFilter) to cover the two new rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C55")) | Out-Null

# Match the author's final selection.
$ws.Range("C40").Select() | Out-Null
